# DeveloperGuide deck maintenance edit:
#  1. Refresh the cached "datetimeFigureOut" footer field text (7/20/17 -> 4/16/2018)
#     on the slide master, every slide layout, and the notes master.
#  2. Remove the now-obsolete "UndoRedoStack" mini-diagram (3 shapes: the blue
#     rectangle "Rectangle 62"/id 59, its outgoing arrow "Straight Arrow
#     Connector 57"/id 61, and the "1" multiplicity label "TextBox 62"/id 63)
#     from slide 1, since the feature moved to VersionedAddressBook.

$p = $ppt.ActivePresentation
$NEW_DATE = "4/16/2018"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh) {
            $sh.TextFrame.TextRange.Text = $NEW_DATE
        }
    }
}

# --- 1a. Slide master's own date footer placeholder ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- 1b. Every slide layout's date footer placeholder ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 1c. Notes master's date footer placeholder ---
try {
    $notesMaster = $p.NotesMaster
    Update-DatePlaceholder $notesMaster.Shapes
} catch {
    Write-Output "NotesMaster update skipped: $_"
}

# --- 2. Remove the obsolete UndoRedoStack mini-diagram shapes from slide 1 ---
$slide = $p.Slides.Item(1)
$idsToRemove = @(59, 61, 63)
foreach ($targetId in $idsToRemove) {
    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            $sh.Delete()
            break
        }
    }
}
